$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "level 2" block: Clay materials (rows 10-14) ---
# Keys pasted into column A first, then values pasted into column B,
# matching the original authoring order (keys column, then values column).
$clayKeys = New-Object 'object[,]' 5,1
$clayKeys[0,0] = "material_clay"
$clayKeys[1,0] = "material_clay_1x2"
$clayKeys[2,0] = "material_clay_1x3"
$clayKeys[3,0] = "material_clay_3x1"
$clayKeys[4,0] = "material_clay_6x1"
$ws.Range("A10:A14").Value = $clayKeys

$clayVals = New-Object 'object[,]' 5,1
$clayVals[0,0] = "Clay"
$clayVals[1,0] = "Clay 1x2"
$clayVals[2,0] = "Clay 1x3"
$clayVals[3,0] = "Clay 3x1"
$clayVals[4,0] = "Clay 6x1"
$ws.Range("B10:B14").Value = $clayVals

# --- "level 1" tweaks: tag_wide / tag_tall (rows 6-7) ---
$wideTallKeys = New-Object 'object[,]' 2,1
$wideTallKeys[0,0] = "tag_wide"
$wideTallKeys[1,0] = "tag_tall"
$ws.Range("A6:A7").Value = $wideTallKeys

$wideTallVals = New-Object 'object[,]' 2,1
$wideTallVals[0,0] = "Wide"
$wideTallVals[1,0] = "Tall"
$ws.Range("B6:B7").Value = $wideTallVals

# --- Iron Block material (row 15) ---
$ws.Range("A15").Value = "material_iron_block"
$ws.Range("B15").Value = "Iron Block"

# --- tag_heavy / tag_light (rows 8-9), entered row by row ---
$ws.Range("A8").Value = "tag_heavy"
$ws.Range("B8").Value = "Heavy"
$ws.Range("A9").Value = "tag_light"
$ws.Range("B9").Value = "Light"

# Final selection left on B16 (the next empty row below the new data)
$ws.Range("B16").Select()
